$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# The query creator now forms the delete/alter-sequence query block starting
# right after row 16, so the three blank placeholder rows that used to sit
# between row 16 and the old row 21 are removed - shifting everything below
# up by three rows (21 -> 18, 22 -> 19, ... 29 -> 26).
$ws.Range("A17:A19").EntireRow.Delete()

# Highlight the generated query rows in bold (new bold Arial font/style).
$ws.Range("A3:K3").Font.Bold = $true
$ws.Range("A8:K8").Font.Bold = $true
$ws.Range("A13:K13").Font.Bold = $true
$ws.Range("A20").Font.Bold = $true
$ws.Range("A25").Font.Bold = $true

# Re-zoom both sheets and move the selection to where the query creator left
# the cursor after forming the new query block.
$ws2.Activate()
$excel.ActiveWindow.Zoom = 130

$ws.Activate()
$excel.ActiveWindow.Zoom = 130
$null = $ws.Range("C26").Select()
